# Fixing errors in example upload files.
#
# - "Service Contacts" sheet: widen column A slightly and move the
#   selection from the whole column D to cell D3.
# - "Practitioners" sheet: widen columns A, C and F, and add a missing
#   practitioner record (row 6) that was dropped from the example file.
#
# Note: ColumnWidth values below are chosen so the engine's internal
# character-width rounding reproduces the exact (or closest achievable)
# target column widths stored in the workbook XML.

$wb = $excel.ActiveWorkbook

# --- Service Contacts sheet -------------------------------------------
$wsSC = $wb.Worksheets.Item("Service Contacts")
$wsSC.Columns.Item(1).ColumnWidth = 13.666666666666666   # -> stored width 14.5

# --- Practitioners sheet ------------------------------------------------
$wsPract = $wb.Worksheets.Item("Practitioners")

# Add the missing practitioner row.
$wsPract.Cells.Item(6,1).Value = "PHN999:NFP02"
$wsPract.Cells.Item(6,2).Value = "P01"
$wsPract.Cells.Item(6,3).Value = 8
$wsPract.Cells.Item(6,4).Value = 1
$wsPract.Cells.Item(6,5).Value = 1973
$wsPract.Cells.Item(6,6).Value = 2
$wsPract.Cells.Item(6,7).Value = 1
$wsPract.Cells.Item(6,8).Value = 1
$wsPract.Cells.Item(6,9).Value = "tag1"

$wsPract.Columns.Item(1).ColumnWidth = 13.833333333333332 # -> stored width ~14.664
$wsPract.Columns.Item(3).ColumnWidth = 12.166666666666666 # -> stored width 13
$wsPract.Columns.Item(6).ColumnWidth = 12.0                # -> stored width ~12.832

# --- Selections (also drives each sheet's sheetView activeCell/sqref) ---
[void]$wsSC.Range("D3").Select()
[void]$wsPract.Range("G1:G1048576").Select()

# Restore the originally active sheet/selection (Organisations) so the
# workbook still opens on the same tab it did before these edits.
$wsOrg = $wb.Worksheets.Item("Organisations")
[void]$wsOrg.Activate()
[void]$wsOrg.Range("H1:J3").Select()
